$d = $word.ActiveDocument

# 1) Update the "Curso (semestre ideal)" line to add the EQD course.
$d.Content.Find.Execute("Curso (semestre ideal): EQN (3)", $false, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EQD (3), EQN (3)", 2)

# 2) Remove the trailing "Requisitos" heading and its "LOQ4073 ..." bullet
#    paragraph (the last two paragraphs of the body, right before the
#    section properties).
$count = $d.Paragraphs.Count
$reqPara = $d.Paragraphs.Item($count - 1)
$r = $d.Range($reqPara.Range.Start, $d.Content.End)
$r.Delete()
